$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values in row 2 (D2, F2, H2)
$ws.Range("D2").Value = 4
$ws.Range("F2").Value = -3
$ws.Range("H2").Value = 46

# Update the active selection to C2
$ws.Range("C2").Select()
